$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "67.528.46"
$ws.Range("E2").Value = "  +0.69%  "

Set-TextValue $ws.Range("D3") "3.855.12"
$ws.Range("E3").Value = "  +0.10%  "

Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.19%  "

Set-TextValue $ws.Range("D5") "456.75"
$ws.Range("E5").Value = "  +7.80%  "

Set-TextValue $ws.Range("D6") "146.25"
$ws.Range("E6").Value = "  +12.68%  "

Set-TextValue $ws.Range("D7") "0.625"
$ws.Range("E7").Value = "  +2.65%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("E9").Value = "  +3.13%  "

$ws.Range("E10").Value = "  -2.93%  "

Set-TextValue $ws.Range("D11") "0.0000318"
$ws.Range("E11").Value = "  -8.72%  "

Set-TextValue $ws.Range("D12") "43.99"
$ws.Range("E12").Value = "  +7.71%  "

Set-TextValue $ws.Range("D13") "10.38"
$ws.Range("E13").Value = "  +2.07%  "

Set-TextValue $ws.Range("D14") "4.465.69"
$ws.Range("E14").Value = "  +0.16%  "

Set-TextValue $ws.Range("D15") "14.85"
$ws.Range("E15").Value = "  -5.66%  "

Set-TextValue $ws.Range("D16") "3.880.00"
$ws.Range("E16").Value = "  +0.95%  "

$ws.Range("E17").Value = "  -0.17%  "

Set-TextValue $ws.Range("D18") "20.07"
$ws.Range("E18").Value = "  +1.55%  "

Set-TextValue $ws.Range("D19") "1.16"
$ws.Range("E19").Value = "  +7.11%  "

Set-TextValue $ws.Range("D20") "67.609.42"
$ws.Range("E20").Value = "  +0.12%  "

Set-TextValue $ws.Range("D21") "427.41"
$ws.Range("E21").Value = "  +4.35%  "

Set-TextValue $ws.Range("D22") "14.82"
$ws.Range("E22").Value = "  -1.88%  "

$ws.Range("E23").Value = "  +6.80%  "

Set-TextValue $ws.Range("D24") "86.90"
$ws.Range("E24").Value = "  +3.06%  "

$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D25") "3.51"
$ws.Range("E25").Value = "  +8.81%  "

$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D26") "10.38"
$ws.Range("E26").Value = "  +11.62%  "

Set-TextValue $ws.Range("D27") "37.67"
$ws.Range("E27").Value = "  +1.18%  "

Set-TextValue $ws.Range("D28") "9.98"
$ws.Range("E28").Value = "  -1.11%  "

$ws.Range("E29").Value = "  +0.91%  "

Set-TextValue $ws.Range("D30") "742.12"
$ws.Range("E30").Value = "  -0.11%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D31") "0.135"
$ws.Range("E31").Value = "  +11.96%  "

$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D32") "13.83"
$ws.Range("E32").Value = "  +7.31%  "

$ws.Range("E33").Value = "  -1.29%  "

Set-TextValue $ws.Range("D34") "42.93"
$ws.Range("E34").Value = "  +10.80%  "

$ws.Range("E35").Value = "  +7.74%  "

Set-TextValue $ws.Range("D36") "57.48"
$ws.Range("E36").Value = "  +3.88%  "

$ws.Range("E37").Value = "  +1.63%  "

$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("E39").Value = "  +4.67%  "

Set-TextValue $ws.Range("D40") "0.356"
$ws.Range("E40").Value = "  +13.71%  "

Set-TextValue $ws.Range("D41") "2.99"
$ws.Range("E41").Value = "  +3.39%  "

Set-TextValue $ws.Range("D42") "2.66"
$ws.Range("E42").Value = "  +17.15%  "

Set-TextValue $ws.Range("D43") "0.0₃0685"
$ws.Range("E43").Value = "  -9.40%  "

$ws.Range("E44").Value = "  +4.92%  "

$ws.Range("E45").Value = "  -0.17%  "

Set-TextValue $ws.Range("D46") "3.45"
$ws.Range("E46").Value = "  +2.33%  "

Set-TextValue $ws.Range("D47") "3.27"
$ws.Range("E47").Value = "  +4.72%  "

$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D48") "2.74"
$ws.Range("E48").Value = "  +7.48%  "

$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D49") "2.14"
$ws.Range("E49").Value = "  +4.20%  "

Set-TextValue $ws.Range("D50") "145.06"
$ws.Range("E50").Value = "  +1.36%  "

Set-TextValue $ws.Range("D51") "2.89"
